# Updates cryptos price list data (Price / Volume(1h) columns, and a couple of
# coin name/link swaps) to match the latest scrape, per commit:
# "Updated cryptos list on Mon Aug 14 10:55:15 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.344.63"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3
$ws.Range("D3").Value = "1.844.47"
$ws.Range("E3").Value = "  -0.18%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9974"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.87"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6269"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9985"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07495"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.62%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2897"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.19%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.34%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07729"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.09%  "

# Row 12
$ws.Range("D12").Value = "1.844.32"
$ws.Range("E12").Value = "  -2.32%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.991"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.73%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6805"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.13%  "

# Row 15
$ws.Range("E15").Value = "  -0.48%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.14%  "

# Row 17
$ws.Range("E17").Value = "  +0.76%  "

# Row 18
$ws.Range("D18").Value = "29.384.30"
$ws.Range("E18").Value = "  -0.07%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.63%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.29%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9984"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.04%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.491"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.28%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9983"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.09%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.427"
$ws.Range("D25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "  -0.93%  "

# Row 27
$ws.Range("E27").Value = "  -0.80%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06484"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +15.68%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.408"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.38%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.476"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.19%  "

# Row 31
$ws.Range("E31").Value = "  -0.36%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.089"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.48%  "

# Row 33
$ws.Range("E33").Value = "  -0.08%  "

# Row 34
$ws.Range("E34").Value = "  -1.86%  "

# Row 35
$ws.Range("E35").Value = "  +0.20%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.578"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.36%  "

# Row 37
$ws.Range("D37").Value = "1.267.88"
$ws.Range("E37").Value = "  +3.36%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.838"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.45%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01835"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.80%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.766"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.29%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9157"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.64%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9979"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.13%  "

# Row 43
$ws.Range("D43").Value = "2.009.34"
$ws.Range("E43").Value = "  -18.35%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.13%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.06%  "

# Row 46
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.079"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.59%  "

# Row 47
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.725"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.29%  "

# Row 48
$ws.Range("E48").Value = "  +1.80%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.3964"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.71%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.982"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.27%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000114"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.37%  "
